# The deck ships with a custom "Integral" theme (ppt/theme/theme1.xml, the
# theme backing the slide master) and the stock "Office Theme"
# (ppt/theme/theme2.xml, backing the notes master). The authored edit swaps
# the two designs: the slide master picks up the default Office Theme
# palette, while the notes master picks up the Integral palette.
#
# PowerPoint's object model surfaces the deck-wide DrawingML colour scheme as
# a 12-slot ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# hanging off a master's Theme. Re-pointing every slot on the slide master's
# theme to the Office Theme's RGB values reproduces that half of the swap.

function Hex2Rgb($hex) {
    # ColorFormat.RGB / ThemeColor.RGB take a COLORREF (0x00BBGGRR), i.e. the
    # byte order is reversed from a "RRGGBB" hex string.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$tcs = $master.Theme.ThemeColorScheme

# Office Theme colour scheme, in ThemeColorScheme slot order.
$officeTheme = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $tcs.Item($i + 1).RGB = Hex2Rgb($officeTheme[$i])
}
